# Add a new "Sheet2" after the existing "Sheet1" and populate it with the
# TestNG DataProvider sample (Test Case Name / UserName / Password table),
# as described in https://www.toolsqa.com/selenium-webdriver/testng-data-provider-excel/

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Inserting with Before=$null, After=$ws1 places the new sheet right after
# Sheet1 and makes it the active tab (mirrors the workbook.xml activeTab="1").
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "Test Case Name"
$ws2.Range("B1").Value = "UserName"
$ws2.Range("C1").Value = "Password"

# Data rows (fill column by column so shared-string order matches)
$ws2.Range("A2").Value = "User 1"
$ws2.Range("A3").Value = "User 2"
$ws2.Range("B2").Value = "testuser1"
$ws2.Range("B3").Value = "testuser2"
$ws2.Range("C2").Value = "Test@123"
$ws2.Range("C3").Value = "Test@123"

# Formatting: bold + centered header, centered data
$ws2.Range("A1:C1").HorizontalAlignment = -4108
$ws2.Range("A1:C1").Font.Bold = $true
$ws2.Range("A2:C3").HorizontalAlignment = -4108

# Column widths
$ws2.Columns("A:C").ColumnWidth = 16.33

# Hyperlink the password cells out to the tutorial site
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://www.toolsqa.com", "", "", "Test@123")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://www.toolsqa.com", "", "", "Test@123")

[void]$ws2.Range("A2").Select()
